$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 7351
$ws.Range("I3").Value = 7489
$ws.Range("J3").Value = 7735
$ws.Range("D4").Value = 1961
$ws.Range("J4").Value = 1681
$ws.Range("J5").Value = 605
$ws.Range("J6").Value = 10561
$ws.Range("D7").Value = 28151
$ws.Range("I7").Value = 26230
$ws.Range("J7").Value = 27933

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J3").Value = 54
$ws.Range("J7").Value = 422

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 464
$ws.Range("J6").Value = 651
$ws.Range("J7").Value = 1761

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J3").Value = 207
$ws.Range("J6").Value = 150
$ws.Range("J7").Value = 560

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J3").Value = 422
$ws.Range("J7").Value = 1267

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J2").Value = 140
$ws.Range("J3").Value = 145
$ws.Range("J7").Value = 400

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J2").Value = 259
$ws.Range("J3").Value = 286
$ws.Range("J6").Value = 250
$ws.Range("J7").Value = 857

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J2").Value = 202
$ws.Range("J7").Value = 703

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J6").Value = 114
$ws.Range("J7").Value = 426

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J2").Value = 221
$ws.Range("J7").Value = 798
$ws.Range("J8").Value = 1761
$ws.Range("J11").Value = 498
$ws.Range("J14").Value = 149
$ws.Range("J18").Value = 226
$ws.Range("J19").Value = 806
$ws.Range("J20").Value = 601
$ws.Range("J23").Value = 256
$ws.Range("I29").Value = 1556
$ws.Range("J29").Value = 1488
$ws.Range("J31").Value = 291
$ws.Range("J33").Value = 1267
$ws.Range("J36").Value = 378
$ws.Range("J37").Value = 857
$ws.Range("J42").Value = 1192
$ws.Range("J44").Value = 220
$ws.Range("J50").Value = 167
$ws.Range("J53").Value = 422
$ws.Range("J54").Value = 552
$ws.Range("J55").Value = 440
$ws.Range("J60").Value = 165
$ws.Range("D63").Value = 345
$ws.Range("J63").Value = 81
$ws.Range("J64").Value = 187
$ws.Range("J65").Value = 703
$ws.Range("J67").Value = 1026
$ws.Range("J72").Value = 107
$ws.Range("J73").Value = 270
$ws.Range("J76").Value = 397
$ws.Range("J78").Value = 323
$ws.Range("J83").Value = 560
$ws.Range("J85").Value = 1145
$ws.Range("J88").Value = 296
$ws.Range("J90").Value = 294
$ws.Range("J91").Value = 320
$ws.Range("J94").Value = 312
$ws.Range("J95").Value = 400
$ws.Range("J99").Value = 426
$ws.Range("D101").Value = 28151
$ws.Range("I101").Value = 26230
$ws.Range("J101").Value = 27933

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J6").Value = 105
$ws.Range("J7").Value = 291

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J6").Value = 283
$ws.Range("J7").Value = 1026

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J3").Value = 111
$ws.Range("J4").Value = 44
$ws.Range("J6").Value = 255
$ws.Range("J7").Value = 552

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I3").Value = 528
$ws.Range("J6").Value = 377
$ws.Range("I7").Value = 1556
$ws.Range("J7").Value = 1488

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J2").Value = 198
$ws.Range("J3").Value = 229
$ws.Range("J6").Value = 312
$ws.Range("J7").Value = 806

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("J4").Value = 14
$ws.Range("J7").Value = 220

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J3").Value = 88
$ws.Range("J7").Value = 397

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("J5").Value = 3
$ws.Range("J6").Value = 61
$ws.Range("J7").Value = 149

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J2").Value = 248
$ws.Range("J3").Value = 240
$ws.Range("J7").Value = 1192

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J6").Value = 101
$ws.Range("J7").Value = 323

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J6").Value = 249
$ws.Range("J7").Value = 440

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("J2").Value = 71
$ws.Range("J7").Value = 256

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J3").Value = 131
$ws.Range("J7").Value = 320

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("J4").Value = 20
$ws.Range("J6").Value = 66
$ws.Range("J7").Value = 187

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J2").Value = 166
$ws.Range("J6").Value = 176
$ws.Range("J7").Value = 601

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("J2").Value = 60
$ws.Range("J7").Value = 226

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J2").Value = 125
$ws.Range("J3").Value = 122
$ws.Range("J7").Value = 378

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J2").Value = 253
$ws.Range("J3").Value = 237
$ws.Range("J7").Value = 798

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J3").Value = 59
$ws.Range("J7").Value = 312

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("J4").Value = 25
$ws.Range("J7").Value = 167

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J4").Value = 29
$ws.Range("J7").Value = 498

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("J3").Value = 67
$ws.Range("J7").Value = 270

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("J3").Value = 54
$ws.Range("J7").Value = 221

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("J6").Value = 157
$ws.Range("J7").Value = 296

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J2").Value = 105
$ws.Range("J7").Value = 294

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("J3").Value = 46
$ws.Range("J7").Value = 165

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J3").Value = 414
$ws.Range("J4").Value = 73
$ws.Range("J7").Value = 1145

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("J6").Value = 40
$ws.Range("J7").Value = 107
